$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data area (A1:G4) first
$ws.Range("A1:G4").Clear()

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 9288.299999999999, 10036, 8426, 0.2073421398798625),
    @(1, 9425.533333333333, 10063, 8369, 0.1745397726694743),
    @(2, 9568, 10247, 8824, 0.1927805582682292),
    @(3, 9873.133333333333, 10544, 8754, 0.1464988231658935),
    @(4, 9114.1, 10237, 7908, 0.1527723471323649),
    @(5, 10040.93333333333, 10948, 9161, 0.167022705078125),
    @(6, 9269.666666666666, 9893, 8459, 0.1634728113810221),
    @(7, 9117.966666666667, 10194, 8317, 0.2142864465713501),
    @(8, 8664.933333333332, 9498, 7300, 0.1868835846583048),
    @(9, 8893.633333333333, 9790, 8124, 0.1887647946675619)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
